$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.523.69"
$ws.Range("E2").Value = "  -3.50%  "

$ws.Range("D3").Value = "3.001.46"
$ws.Range("E3").Value = "  -2.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.03"
$ws.Range("E5").Value = "  -0.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.18"
$ws.Range("E6").Value = "  -1.11%  "

$ws.Range("D8").Value = "2.999.48"
$ws.Range("E8").Value = "  -2.72%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.496"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.149"
$ws.Range("E10").Value = "  -4.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.14"
$ws.Range("E11").Value = "  -1.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("E12").Value = "  -2.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000223"
$ws.Range("E13").Value = "  -2.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.08"
$ws.Range("E14").Value = "  -2.20%  "

$ws.Range("D15").Value = "3.483.95"
$ws.Range("E15").Value = "  -2.89%  "

$ws.Range("E16").Value = "  -1.40%  "

$ws.Range("D17").Value = "61.533.19"
$ws.Range("E17").Value = "  -3.47%  "

$ws.Range("D18").Value = "2.993.11"
$ws.Range("E18").Value = "  -3.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.64"
$ws.Range("E19").Value = "  -0.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.45"
$ws.Range("E20").Value = "  -4.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.25"
$ws.Range("E21").Value = "  -1.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.678"
$ws.Range("E22").Value = "  -3.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.96"
$ws.Range("E23").Value = "  -3.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.29"
$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.01"
$ws.Range("E25").Value = "  -2.20%  "

$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("E27").Value = "  -1.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.81"
$ws.Range("E28").Value = "  -6.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.89"
$ws.Range("E30").Value = "  -1.23%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.64"
$ws.Range("E31").Value = "  -2.53%  "

$ws.Range("B32").Value = "Mantle"
$ws.Range("C32").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  +2.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.54"
$ws.Range("E33").Value = "  +1.17%  "

$ws.Range("E34").Value = "  -5.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "54.78"
$ws.Range("E35").Value = "  -4.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.91"
$ws.Range("E36").Value = "  -2.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "453.99"
$ws.Range("E37").Value = "  -8.12%  "

$ws.Range("D38").Value = "3.170.65"
$ws.Range("E38").Value = "  -3.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0789"
$ws.Range("E39").Value = "  -1.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.121"
$ws.Range("E40").Value = "  +2.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0386"
$ws.Range("E41").Value = "  -3.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.15"
$ws.Range("E42").Value = "  -0.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.47"
$ws.Range("E43").Value = "  -6.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.14"
$ws.Range("E44").Value = "  +8.79%  "

$ws.Range("E45").Value = "  +0.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.245"
$ws.Range("E46").Value = "  -5.15%  "

$ws.Range("E47").Value = "  -3.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.77"
$ws.Range("E48").Value = "  -2.37%  "

$ws.Range("E49").Value = "  -1.04%  "

$ws.Range("D50").Value = "0.0₃0498"
$ws.Range("E50").Value = "  -7.80%  "

$ws.Range("E51").Value = "  +6.00%  "
